$wb = $excel.ActiveWorkbook

# --- Add new "Player Info" sheet, positioned before "ODI Batting" ---
$odiBattingRef = $wb.Worksheets.Item("ODI Batting")
$newSheet = $wb.Worksheets.Add($odiBattingRef)
$newSheet.Name = "Player Info"

# Re-fetch sheet references fresh by name (adding a sheet invalidates old refs)
$playerInfo = $wb.Worksheets.Item("Player Info")
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBowling = $wb.Worksheets.Item("ODI Bowling")

# A guaranteed-empty cell on each sheet, used as a "blank formatting" source
# so that Range.Copy(...) can reset a cell's style to the default (no
# explicit style) without touching the text value already stored there.
$blankBatting = $odiBatting.Range("ZZ999")
$blankBowling = $odiBowling.Range("ZZ999")
$blankPlayerInfo = $playerInfo.Range("ZZ999")

# Copy header row (A1:D1) from "ODI Batting" into "Player Info" to inherit the
# bold/border/centered header style used throughout the workbook, then
# overwrite with the new header text (style is kept, value is replaced).
$odiBatting.Range("A1:D1").Copy($playerInfo.Range("A1:D1"))
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row for "Player Info" (plain/default style, text values)
$playerInfo.Range("A2").Value = "'4648"
$blankPlayerInfo.Copy($playerInfo.Range("A2"))
$playerInfo.Range("B2").Value = "Ashton James Turner"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# --- Update "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code ---
$odiBatting.Range("D1").Value = "MATCH_CODE"
$battingRows = $odiBatting.UsedRange.Rows.Count
for ($r = 2; $r -le $battingRows; $r++) {
    $cell = $odiBatting.Cells.Item($r, 4)
    $val = $cell.Text
    if ($val -match "MatchCode=(\d+)") {
        $cell.Value = "'" + $matches[1]
        $blankBatting.Copy($cell)
    }
}

# --- Update "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code ---
$odiBowling.Range("B1").Value = "MATCH_CODE"
$bowlingRows = $odiBowling.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingRows; $r++) {
    $cell = $odiBowling.Cells.Item($r, 2)
    $val = $cell.Text
    if ($val -match "MatchCode=(\d+)") {
        $cell.Value = "'" + $matches[1]
        $blankBowling.Copy($cell)
    }
}
